$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("references")
$ws.Range("A2").Value = 'Verify reference "Cleveland Clinic. Allergy Overview."'
$ws.Range("A3").Value = 'Verify reference "American College of Allergy, Asthma, & Immunology. Allergy Symptoms."'
$ws.Range("A4").Value = 'Verify reference "Asthma and Allergy Foundation of America. Pollen Allergy."'
$ws.Range("A5").Value = 'Verify reference "Asthma and Allergy Foundation of America. Pet Allergy"'
$ws.Range("A6").Value = 'Verify reference "American College of Allergy, Asthma, & Immunology."'
$ws.Range("A7").Value = 'Verify reference "Mayo Clinic. Cold or allergy"'
$ws.Activate()
